$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1400.3334
$ws.Range("J19").Value = 1202
$ws.Range("L19").Value = 1202
$ws.Range("N19").Value = -1552
$ws.Range("H33").Value = 975.875
$ws.Range("I33").Value = 555.9091
$ws.Range("K33").Value = 555.9091
$ws.Range("M33").Value = -326.9091
$ws.Range("H42").Value = 473.625
$ws.Range("J42").Value = 400
$ws.Range("L42").Value = 1200
$ws.Range("N42").Value = -1660
$ws.Range("H51").Value = 500000000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H132").Value = 4953.9414
$ws.Range("I132").Value = 5825.5
$ws.Range("K132").Value = 17476.5
$ws.Range("M132").Value = -14946.5
$ws.Range("H135").Value = 996.41174
$ws.Range("I135").Value = 828.8570999999999
$ws.Range("J135").Value = 1778.3334
$ws.Range("K135").Value = 7459.7139
$ws.Range("L135").Value = 16005.0006
$ws.Range("M135").Value = -4924.7139
$ws.Range("N135").Value = -21075.0006
$ws.Range("H138").Value = 4213.825
$ws.Range("I138").Value = 5287.8887
$ws.Range("J138").Value = 3902
$ws.Range("K138").Value = 15863.6661
$ws.Range("L138").Value = 11706
$ws.Range("M138").Value = -10723.6661
$ws.Range("N138").Value = -21986

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2077.0244
$ws.Range("I32").Value = 2137.4614
$ws.Range("K32").Value = 2137.4614
$ws.Range("M32").Value = -1850.4614
$ws.Range("H45").Value = 34323.152
$ws.Range("I45").Value = 48601.223
$ws.Range("J45").Value = 2197.5
$ws.Range("K45").Value = 48601.223
$ws.Range("L45").Value = 2197.5
$ws.Range("M45").Value = -48224.223
$ws.Range("N45").Value = -2951.5
$ws.Range("H61").Value = 3812.724
$ws.Range("I61").Value = 1989.1818
$ws.Range("J61").Value = 4927.1113
$ws.Range("K61").Value = 1989.1818
$ws.Range("L61").Value = 4927.1113
$ws.Range("M61").Value = -1777.1818
$ws.Range("N61").Value = -5351.1113
$ws.Range("H110").Value = 789.2
$ws.Range("I110").Value = 719.3333
$ws.Range("K110").Value = 719.3333
$ws.Range("M110").Value = 1325.6667
$ws.Range("H122").Value = 13581.1
$ws.Range("I122").Value = 19383
$ws.Range("K122").Value = 58149
$ws.Range("M122").Value = -55699
$ws.Range("H132").Value = 5637.5806
$ws.Range("I132").Value = 4442.1177
$ws.Range("J132").Value = 7089.2144
$ws.Range("K132").Value = 13326.3531
$ws.Range("L132").Value = 21267.6432
$ws.Range("M132").Value = -10796.3531
$ws.Range("N132").Value = -26327.6432
$ws.Range("H136").Value = 3812.724
$ws.Range("I136").Value = 1989.1818
$ws.Range("J136").Value = 4927.1113
$ws.Range("K136").Value = 5967.5454
$ws.Range("L136").Value = 14781.3339
$ws.Range("M136").Value = -3417.5454
$ws.Range("N136").Value = -19881.3339

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2665.4707
$ws.Range("I86").Value = 2269.375
$ws.Range("K86").Value = 2269.375
$ws.Range("M86").Value = -1146.375
$ws.Range("H89").Value = 2665.4707
$ws.Range("I89").Value = 2269.375
$ws.Range("K89").Value = 11346.875
$ws.Range("M89").Value = -5730.875
$ws.Range("H105").Value = 9631620
$ws.Range("I105").Value = 668483.9
$ws.Range("K105").Value = 668483.9
$ws.Range("M105").Value = -666736.9

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 19250
$ws.Range("H58").Value = 3879.6
$ws.Range("I58").Value = 3398
$ws.Range("K58").Value = 3398
$ws.Range("M58").Value = -3195
$ws.Range("H99").Value = 1749.75
$ws.Range("I99").Value = 1700
$ws.Range("J99").Value = 1799.5
$ws.Range("K99").Value = 1700
$ws.Range("L99").Value = 1799.5
$ws.Range("M99").Value = -202
$ws.Range("N99").Value = -4795.5
$ws.Range("H107").Value = 2273429
$ws.Range("I107").Value = 4167079.5
$ws.Range("J107").Value = 1048.2
$ws.Range("K107").Value = 4167079.5
$ws.Range("L107").Value = 1048.2
$ws.Range("M107").Value = -4165159.5
$ws.Range("N107").Value = -4888.2
$ws.Range("H120").Value = 34499.5
$ws.Range("J120").Value = 34499.5
$ws.Range("L120").Value = 34499.5
$ws.Range("N120").Value = -41757.5
$ws.Range("H121").Value = 32499
$ws.Range("J121").Value = 32499
$ws.Range("L121").Value = 32499
$ws.Range("N121").Value = -35119
$ws.Range("H126").Value = 1749.75
$ws.Range("I126").Value = 1700
$ws.Range("J126").Value = 1799.5
$ws.Range("K126").Value = 5100
$ws.Range("L126").Value = 5398.5
$ws.Range("M126").Value = -2630
$ws.Range("N126").Value = -10338.5
$ws.Range("H136").Value = 3879.6
$ws.Range("I136").Value = 3398
$ws.Range("K136").Value = 10194
$ws.Range("M136").Value = -7644
$ws.Range("H141").Value = 192766.4
$ws.Range("J141").Value = 196535.42
$ws.Range("L141").Value = 196535.42
$ws.Range("N141").Value = -206895.42

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 97225360
$ws.Range("J32").Value = 16670282
$ws.Range("L32").Value = 50010846
$ws.Range("N32").Value = -50011412
$ws.Range("H37").Value = 2188746.2
$ws.Range("J37").Value = 2188746.2
$ws.Range("L37").Value = 6566238.600000001
$ws.Range("N37").Value = -6566462.600000001
$ws.Range("H51").Value = 3457.0833
$ws.Range("I51").Value = 246.25
$ws.Range("K51").Value = 738.75
$ws.Range("M51").Value = -278.75
$ws.Range("H59").Value = 4410
$ws.Range("I59").Value = 845
$ws.Range("J59").Value = 5428.5713
$ws.Range("K59").Value = 2535
$ws.Range("L59").Value = 16285.7139
$ws.Range("M59").Value = -1995
$ws.Range("N59").Value = -17365.7139
$ws.Range("H132").Value = 2877.6
$ws.Range("I132").Value = 1695
$ws.Range("J132").Value = 3666
$ws.Range("K132").Value = 15255
$ws.Range("L132").Value = 32994
$ws.Range("M132").Value = -12725
$ws.Range("N132").Value = -38054
$ws.Range("H138").Value = 6101.8335
$ws.Range("I138").Value = 6101.8335
$ws.Range("K138").Value = 18305.5005
$ws.Range("M138").Value = -13165.5005

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 35716776
$ws.Range("I80").Value = 100001580
$ws.Range("J80").Value = 2992.5
$ws.Range("K80").Value = 100001580
$ws.Range("L80").Value = 2992.5
$ws.Range("M80").Value = -100000582
$ws.Range("N80").Value = -4988.5
$ws.Range("H83").Value = 35716776
$ws.Range("I83").Value = 100001580
$ws.Range("J83").Value = 2992.5
$ws.Range("K83").Value = 500007900
$ws.Range("L83").Value = 14962.5
$ws.Range("M83").Value = -500002908
$ws.Range("N83").Value = -24946.5
$ws.Range("H102").Value = 15653.5
$ws.Range("I102").Value = 3380
$ws.Range("J102").Value = 36109.332
$ws.Range("K102").Value = 3380
$ws.Range("L102").Value = 36109.332
$ws.Range("M102").Value = -1758
$ws.Range("N102").Value = -39353.332
$ws.Range("H132").Value = 2560.9546
$ws.Range("I132").Value = 2021.5
$ws.Range("K132").Value = 6064.5
$ws.Range("M132").Value = -3534.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2249.6667
$ws.Range("I46").Value = 2374.5
$ws.Range("K46").Value = 2374.5
$ws.Range("M46").Value = -2186.5
$ws.Range("H55").Value = 560.36365
$ws.Range("I55").Value = 316.42856
$ws.Range("J55").Value = 987.25
$ws.Range("K55").Value = 316.42856
$ws.Range("L55").Value = 987.25
$ws.Range("M55").Value = -143.42856
$ws.Range("N55").Value = -1333.25
$ws.Range("H132").Value = 7995.5884
$ws.Range("I132").Value = 5326.4443
$ws.Range("J132").Value = 10998.375
$ws.Range("K132").Value = 15979.3329
$ws.Range("L132").Value = 32995.125
$ws.Range("M132").Value = -13449.3329
$ws.Range("N132").Value = -38055.125

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 35714804
$ws.Range("J100").Value = 111111816
$ws.Range("L100").Value = 222223632
$ws.Range("N100").Value = -222224714
$ws.Range("H132").Value = 2989.8333
$ws.Range("I132").Value = 2478.9524
$ws.Range("K132").Value = 7436.8572
$ws.Range("M132").Value = -4906.8572
$ws.Range("H136").Value = 18833.375
$ws.Range("I136").Value = 20690.379
$ws.Range("J136").Value = 15999
$ws.Range("K136").Value = 62071.137
$ws.Range("L136").Value = 47997
$ws.Range("M136").Value = -59521.137
$ws.Range("N136").Value = -53097
$ws.Range("H138").Value = 113499.25
$ws.Range("J138").Value = 117999
$ws.Range("L138").Value = 117999
$ws.Range("N138").Value = -128279
